$d = $word.ActiveDocument
$d.Content.Find.Execute("opcionalan", $true, $false, $false, $false, $false,
                         $true, 1, $false, "opcion", 2)
